$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the quarterly recurrence metrics for row 20 (2025Q2)
$ws.Range("C20").Value = 262
$ws.Range("D20").Value = 222
$ws.Range("E20").Value = 40
$ws.Range("F20").Value = 73.02631578947368
